# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# This script updates the "Metadata" sheet (sheet1) values and completely
# rebuilds the "Include #0" sheet (sheet2) content - converting it from a
# SNOMED CT "is-a" filter definition into an explicit concept enumeration
# of the Fitzpatrick skin type codes.
#
# Both sheets' contents are cleared and then rewritten cell-by-cell, in
# strict reading order (top-to-bottom, left-to-right, sheet1 before
# sheet2). The underlying engine builds its shared-string table in
# first-use order, so writing every cell back (even the ones whose text
# did not change) reproduces the same shared-string ordering as a clean
# rebuild of the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Metadata
$ws2 = $wb.Worksheets.Item(2)   # Include #0

# xlPasteFormats
$xlPasteFormats = -4122

# Clear all existing values (formatting/styles on each cell are kept).
$ws1.Range("A1:B14").ClearContents()
$ws2.Range("A1:C4").ClearContents()

# Drop the now-unused third column on sheet2 entirely.
$ws2.Columns.Item(3).Delete()

# Sheet2 needs 5 additional rows (5-9). Give them the same body style
# (s="2") as the rest of the table before any values are written.
$bodyFormat = $ws2.Range("A2:B2")
foreach ($r in 5..9) {
    $bodyFormat.Copy()
    $ws2.Range("A" + $r + ":B" + $r).PasteSpecial($xlPasteFormats)
}

# -----------------------------------------------------------------
# Sheet 1 ("Metadata")
# -----------------------------------------------------------------
$ws1.Range("A1").Value = "Property"
$ws1.Range("B1").Value = "Value"

$ws1.Range("A2").Value = "URL"
$ws1.Range("B2").Value = "https://clinyqai.github.io/open-nursing-core-ig/ValueSet/skintone-vs"

$ws1.Range("A3").Value = "Version"
$ws1.Range("B3").Value = "0.1.0"

$ws1.Range("A4").Value = "Name"
$ws1.Range("B4").Value = "SkinToneVS"

$ws1.Range("A5").Value = "Title"
$ws1.Range("B5").Value = "Fitzpatrick Skin Tone Value Set"

$ws1.Range("A6").Value = "Status"
$ws1.Range("B6").Value = "draft"

$ws1.Range("A7").Value = "Experimental"
# "false" would otherwise be auto-interpreted as a Boolean by the Excel
# engine (producing t="b") instead of being stored as plain text like the
# target file expects (t="s"). Force it in as text via a leading
# apostrophe, then re-apply the plain body formatting (copied from a
# neighboring cell) so the quote-prefix flag doesn't linger on the style.
$ws1.Range("B7").Value = "'false"
$ws1.Range("A6").Copy()
$ws1.Range("B7").PasteSpecial($xlPasteFormats)

$ws1.Range("A8").Value = "Date"
$ws1.Range("B8").Value = "2025-12-26T14:13:58+00:00"

$ws1.Range("A9").Value = "Publisher"
# B9 stays blank.

$ws1.Range("A10").Value = "Jurisdiction"
# B10 is an (empty-string) shared-string cell rather than a truly blank
# cell. A lone leading apostrophe produces an empty text value that
# reuses the shared empty string; re-apply the plain body formatting
# afterwards for the same reason as above.
$ws1.Range("B10").Value = "'"
$ws1.Range("A9").Copy()
$ws1.Range("B10").PasteSpecial($xlPasteFormats)

$ws1.Range("A11").Value = "Description"
$ws1.Range("B11").Value = "Value set for Fitzpatrick skin type classifications"

$ws1.Range("A12").Value = "Purpose"
# B12 stays blank.

$ws1.Range("A13").Value = "Copyright"
# B13 stays blank.

$ws1.Range("A14").Value = "Immutable"
$ws1.Range("B14").Value = "BooleanType[null]"

# -----------------------------------------------------------------
# Sheet 2 ("Include #0")
# -----------------------------------------------------------------
$ws2.Range("A1").Value = "Concept"
$ws2.Range("B1").Value = "Description"

$ws2.Range("A2").Value = "fitzpatrick-1"
$ws2.Range("B2").Value = "Type I"

$ws2.Range("A3").Value = "fitzpatrick-2"
$ws2.Range("B3").Value = "Type II"

$ws2.Range("A4").Value = "fitzpatrick-3"
$ws2.Range("B4").Value = "Type III"

$ws2.Range("A5").Value = "fitzpatrick-4"
$ws2.Range("B5").Value = "Type IV"

$ws2.Range("A6").Value = "fitzpatrick-5"
$ws2.Range("B6").Value = "Type V"

$ws2.Range("A7").Value = "fitzpatrick-6"
$ws2.Range("B7").Value = "Type VI"

# Blank separator row (row 8), stored as references to the shared empty
# string (matching the original blank separator row in this sheet).
$ws2.Range("A8").Value = "'"
$ws2.Range("B8").Value = "'"
$ws2.Range("A2").Copy()
$ws2.Range("A8").PasteSpecial($xlPasteFormats)
$ws2.Range("A2").Copy()
$ws2.Range("B8").PasteSpecial($xlPasteFormats)

$ws2.Range("A9").Value = "System URI"
$ws2.Range("B9").Value = "https://clinyqai.github.io/open-nursing-core-ig/CodeSystem/onc-observation-codes"
